$d = $word.ActiveDocument

# --- Change 1: Professional summary paragraph ---
# "affecting all Black and Asian-American voters," -> "affecting 50M voters,"
$p4 = $d.Paragraphs.Item(4).Range
$ok1 = $p4.Find.Execute("affecting all Black and Asian-American voters,", $true, $false, $false, $false, $false, $true, 1, $false, "affecting 50M voters,", 2)
Write-Output "Summary paragraph updated: $ok1"

# --- Change 2: Work experience bullet (requires bold/colored "50M" run) ---
# "... affecting all Black and Asian-American voters, developed ..." ->
# "... affecting " + bold/colored "50M" + " voters, developed ..."
$p10 = $d.Paragraphs.Item(10).Range
$rngBullet = $p10.Duplicate
$ok2 = $rngBullet.Find.Execute("all Black and Asian-American", $true, $false, $false, $false, $false, $true, 1, $false, "50M", 2)
Write-Output "Bullet text replaced: $ok2"
$rngBullet.Font.Bold = -1
$rngBullet.Font.Color = 5258796

# --- Change 3: Key project impact statement ---
# "affecting all Black and Asian-American voters," -> "affecting 50M voters nationwide,"
$p47 = $d.Paragraphs.Item(47).Range
$ok3 = $p47.Find.Execute("affecting all Black and Asian-American voters,", $true, $false, $false, $false, $false, $true, 1, $false, "affecting 50M voters nationwide,", 2)
Write-Output "Impact statement updated: $ok3"
